$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number + date range)
$ws.Range("A8").Value = "Volume 31   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/24/2024  Through  6/30/2024"

# Data table updates (rows 14-30)
$ws.Range("N14").Value = -60
$ws.Range("D15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 18
$ws.Range("K15").Value = -14.285714285714
$ws.Range("L15").Value = 12.5
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -61.702127659574
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 5
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -48
$ws.Range("I16").Value = 142
$ws.Range("J16").Value = 133
$ws.Range("K16").Value = 6.766917293233
$ws.Range("L16").Value = -19.318181818181
$ws.Range("M16").Value = -29.702970297029
$ws.Range("N16").Value = -88.166666666666
$ws.Range("C17").Value = 23
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 130
$ws.Range("F17").Value = 70
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = 37.254901960784
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 328
$ws.Range("K17").Value = 21.951219512195
$ws.Range("L17").Value = 30.718954248366
$ws.Range("M17").Value = 82.648401826484
$ws.Range("N17").Value = -33.444259567387
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 107
$ws.Range("J18").Value = 105
$ws.Range("K18").Value = 1.904761904761
$ws.Range("L18").Value = 1.904761904761
$ws.Range("M18").Value = -41.530054644808
$ws.Range("N18").Value = -89.489194499017
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -22.033898305084
$ws.Range("I19").Value = 292
$ws.Range("J19").Value = 348
$ws.Range("K19").Value = -16.091954022988
$ws.Range("L19").Value = -17.514124293785
$ws.Range("M19").Value = 9.774436090225
$ws.Range("N19").Value = -20.218579234972
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 14.285714285714
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 36.842105263157
$ws.Range("I20").Value = 111
$ws.Range("J20").Value = 126
$ws.Range("K20").Value = -11.904761904761
$ws.Range("L20").Value = -23.972602739726
$ws.Range("M20").Value = -22.377622377622
$ws.Range("N20").Value = -89.488636363636
$ws.Range("C21").Value = 46
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 27.777777777777
$ws.Range("F21").Value = 174
$ws.Range("G21").Value = 174
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 1078
$ws.Range("J21").Value = 1069
$ws.Range("K21").Value = 0.841908325537
$ws.Range("L21").Value = -3.057553956834
$ws.Range("M21").Value = 3.653846153846
$ws.Range("N21").Value = -74.976787372330
$ws.Range("D22").Value = 1
$ws.Range("G14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("H14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("G14").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = -100
$ws.Range("H14").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = -20
$ws.Range("D23").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("G23").Value = 1
$ws.Range("L23").Value = -66.666666666666
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -20
$ws.Range("F24").Value = 110
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = 35.802469135802
$ws.Range("I24").Value = 631
$ws.Range("J24").Value = 606
$ws.Range("K24").Value = 4.125412541254
$ws.Range("L24").Value = 32.285115303983
$ws.Range("M24").Value = 33.686440677966
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 6
$ws.Range("G14").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = -16.666666666666
$ws.Range("H14").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 91
$ws.Range("J25").Value = 117
$ws.Range("K25").Value = -22.222222222222
$ws.Range("L25").Value = -4.210526315789
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 19
$ws.Range("E26").Value = 10.526315789473
$ws.Range("F26").Value = 75
$ws.Range("G26").Value = 64
$ws.Range("H26").Value = 17.1875
$ws.Range("I26").Value = 452
$ws.Range("J26").Value = 378
$ws.Range("K26").Value = 19.576719576719
$ws.Range("L26").Value = 21.505376344086
$ws.Range("M26").Value = -1.952277657266
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 7
$ws.Range("H27").Value = 16.666666666666
$ws.Range("I27").Value = 29
$ws.Range("K27").Value = 7.407407407407
$ws.Range("L27").Value = 26.086956521739
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 4
$ws.Range("G14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -50
$ws.Range("H14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 42.857142857142
$ws.Range("I28").Value = 41
$ws.Range("J28").Value = 39
$ws.Range("K28").Value = 5.128205128205
$ws.Range("L28").Value = 0
$ws.Range("C29").Value = 3
$ws.Range("G14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("F29").Value = 8
$ws.Range("H29").Value = 166.666666666667
$ws.Range("I29").Value = 20
$ws.Range("K29").Value = -13.043478260869
$ws.Range("L29").Value = -4.761904761904
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -79.381443298969
$ws.Range("C30").Value = 1
$ws.Range("G14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("I30").Value = 17
$ws.Range("K30").Value = -10.526315789473
$ws.Range("L30").Value = 6.25
$ws.Range("M30").Value = -34.615384615384
$ws.Range("N30").Value = -81.521739130434
